$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells remain text so values like "298.39" and "0.50%"
# are stored as inline strings (matching the original text-typed cells),
# not auto-converted to numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.50%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.25%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.118"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.11%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08039"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.58%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.570"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "53.51%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.60%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.823"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.59%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9177"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.01%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1731"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.97%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07316"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08329"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.32%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03033"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.31%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09964"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.58%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001491"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005912"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.85%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.504"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.62%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.251"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.25%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.35%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.78%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.631"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.70%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04567"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.56%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001258"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.50%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.53%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-9.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003430"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "83.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01839"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.84%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04515"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.14%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007030"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.22%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.94%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.45%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009815"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.45%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006490"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.77%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-56.67%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
